# "Log all coming smss and test project on real condition"
#
# Adds two more logged devices to the "serials" sheet (rows 5 & 6):
#   - row 5: a brand new "Second device" entry (reference 106, CC1234-DD1245,
#            dated 2021-09-11)
#   - row 6: another real-world duplicate of the existing "New 20 device"
#            entry (reference 101, JM200-Jm299, dated 2012-07-02)
# and widens the "invalids" sheet's first column so the longer reference
# strings it lists are fully visible.

$wb = $excel.ActiveWorkbook

$serials = $wb.Worksheets.Item("serials")

# --- Row 5: brand-new device -------------------------------------------
$serials.Range("A5").Value = 5
$serials.Range("B5").Value = 106
$serials.Range("C5").Value = "Second device"
$serials.Range("D5").Value = "CC1234"
$serials.Range("E5").Value = "DD1245"
$serials.Range("F5").Value = 44450

# --- Row 6: another real-condition sample --------------------------------
$serials.Range("A6").Value = 6
$serials.Range("B6").Value = 101
$serials.Range("C6").Value = "New 20 device"
$serials.Range("D6").Value = "JM200"
$serials.Range("E6").Value = "Jm299"
$serials.Range("F6").Value = 41092

# --- Widen "invalids" column A so entries are readable, and update its
#     selection before returning focus (and the active-cell selection) to
#     "serials", which stays the active sheet.
$invalids = $wb.Worksheets.Item("invalids")
$invalids.Activate()
$invalids.Columns.Item(1).ColumnWidth = 24.25
$invalids.Range("A4").Select()

$serials.Activate()
$serials.Range("E6").Select()
